$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows down
$ws.Rows.Item(4).Insert()

# Set the new cell's content (keep as text, not auto-converted to a date serial)
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "02/18/2022"

# Update the selection to D2 (no explicit top-left cell scroll)
$ws.Range("D2").Select()
